# Auto-generated edit script for cs-en-us-007pct.xlsx weekly refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering date range) ---
$ws.Range("A8").Value = "Volume 32   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/5/2025  Through  5/11/2025"

# --- Cells changing type from a text placeholder ("0"/"***.*") to a real number:
# clone the style from a same-format donor cell (row 29, untouched elsewhere in this edit)
# via Copy (which carries over the cell style), then overwrite with the new numeric value. ---
$ws.Range("I29").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 3
$ws.Range("I29").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 3
$ws.Range("I29").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 1
$ws.Range("K29").Copy($ws.Range("E20"))
$ws.Range("E20").Value = 200
$ws.Range("I29").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K29").Copy($ws.Range("E22"))
$ws.Range("E22").Value = 0

# --- Cells changing type from a real number to a text placeholder ("0"/"***.*"):
# clone the whole cell (value + style) from a same-placeholder donor cell (row 29). ---
$ws.Range("C29").Copy($ws.Range("D27"))
$ws.Range("E29").Copy($ws.Range("E27"))
$ws.Range("C29").Copy($ws.Range("D28"))
$ws.Range("E29").Copy($ws.Range("E28"))
$ws.Range("C29").Copy($ws.Range("G31"))
$ws.Range("E29").Copy($ws.Range("H31"))

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -30.76923076923
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 17.5
$ws.Range("L16").Value = 11.904761904761
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -82.720588235294
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -87.5
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -51.851851851851
$ws.Range("J17").Value = 93
$ws.Range("K17").Value = -17.204301075268
$ws.Range("L17").Value = 11.59420289855
$ws.Range("M17").Value = 48.076923076923
$ws.Range("N17").Value = 1.315789473684
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 49
$ws.Range("K18").Value = -18.367346938775
$ws.Range("L18").Value = -14.893617021276
$ws.Range("M18").Value = 90.47619047619
$ws.Range("N18").Value = -67.479674796748
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 9.375
$ws.Range("I19").Value = 132
$ws.Range("J19").Value = 141
$ws.Range("K19").Value = -6.382978723404
$ws.Range("L19").Value = -29.787234042553
$ws.Range("M19").Value = 78.378378378378
$ws.Range("N19").Value = -12
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = -65
$ws.Range("L20").Value = -56.25
$ws.Range("M20").Value = -69.565217391304
$ws.Range("N20").Value = -94.696969696969
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -16
$ws.Range("F21").Value = 70
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -19.540229885057
$ws.Range("I21").Value = 308
$ws.Range("J21").Value = 347
$ws.Range("K21").Value = -11.239193083573
$ws.Range("L21").Value = -16.076294277929
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = -59.420289855072
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 11
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 83.333333333333
$ws.Range("L22").Value = 37.5
$ws.Range("M22").Value = 266.666666666667
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -85.714285714285
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 62
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = -22.5
$ws.Range("L23").Value = 31.914893617021
$ws.Range("M23").Value = 12.727272727272
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 4.545454545454
$ws.Range("F24").Value = 88
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -13.725490196078
$ws.Range("I24").Value = 381
$ws.Range("J24").Value = 426
$ws.Range("K24").Value = -10.56338028169
$ws.Range("L24").Value = 3.814713896457
$ws.Range("M24").Value = 62.127659574468
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 56
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = -20
$ws.Range("I25").Value = 221
$ws.Range("J25").Value = 272
$ws.Range("K25").Value = -18.75
$ws.Range("L25").Value = 13.333333333333
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = -36.842105263157
$ws.Range("I26").Value = 138
$ws.Range("J26").Value = 163
$ws.Range("K26").Value = -15.337423312883
$ws.Range("L26").Value = -0.719424460431
$ws.Range("M26").Value = 27.777777777777
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("L27").Value = -70
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 20
$ws.Range("K28").Value = 53.846153846153
